$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (from row 23 through row 245) is changed from a per-month
# running index (reset to 0 at the start of each month) to a single
# continuous running index across the whole year. For every affected
# row, the new value equals (row number - 2).
for ($r = 23; $r -le 245; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
